$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") on every data row (2-332) is bumped from
# serial date 46074 (2026-02-21) to 46075 (2026-02-22).
for ($r = 2; $r -le 332; $r++) {
    $ws.Cells.Item($r, 3).Value = 46075
}
